# LLF_sym now properly takes hazard function as parameter
# Rename the covariate-data header columns so the hazard-function related
# columns use the new naming convention (T, kVec, cVec, eVec, fVec) instead
# of the old one (time, kVec, cVec/Evec/FVec/Fvec).

$wb = $excel.ActiveWorkbook

# --- DS1 sheet ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("DS1")
$ws1.Range("A1").Value = "T"
$ws1.Range("B1").Value = "kVec"
$ws1.Range("C1").Value = "cVec"
$ws1.Range("D1").Value = "eVec"
$ws1.Range("E1").Value = "fVec"

# Update the active selection on DS1 to match the edited workbook.
$ws1.Range("E1").Select()

# --- DS2 sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DS2")
$ws2.Range("A1").Value = "T"
$ws2.Range("B1").Value = "kVec"
$ws2.Range("C1").Value = "fVec"
$ws2.Range("D1").Value = "eVec"
$ws2.Range("E1").Value = "cVec"

# DS2 is the active sheet/tab; update its selection too.
$ws2.Activate()
$ws2.Range("E2").Select()
